$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.8983344824904379
$ws1.Range("C2").Value = -1.044205741479208
$ws1.Range("B3").Value = 1.41889673321519
$ws1.Range("C3").Value = -0.9251340377201603
$ws1.Range("B4").Value = 0.7289269858143762
$ws1.Range("C4").Value = -0.6697844404295697

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.684346137949199
$ws2.Range("C2").Value = -0.2576497276055958
$ws2.Range("B3").Value = 1.462435036328542
$ws2.Range("C3").Value = 0.6404751458697506
$ws2.Range("B4").Value = -0.3869909651594174
$ws2.Range("C4").Value = -0.7172796495907583
